$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 170
$ws1.Range("F13").Value = 5281
$ws1.Range("F17").Value = 2764
$ws1.Range("F18").Value = 2764
$ws1.Range("F22").Value = 3992
$ws1.Range("F37").Value = 6941
$ws1.Range("F42").Value = 1412
$ws1.Range("F44").Value = 731
$ws1.Range("F46").Value = 2365
$ws1.Range("F47").Value = 318
$ws1.Range("F50").Value = 792
$ws1.Range("F51").Value = 943

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 170
$ws4.Range("F19").Value = 2764
$ws4.Range("F20").Value = 2764
$ws4.Range("F24").Value = 3992
$ws4.Range("F36").Value = 6941
$ws4.Range("F42").Value = 1412
$ws4.Range("F44").Value = 731
$ws4.Range("F46").Value = 2365
$ws4.Range("F47").Value = 318
$ws4.Range("F49").Value = 792
$ws4.Range("F50").Value = 943
